$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "66.683.93"
Set-TextValue "E2" "  +1.30%  "

Set-TextValue "D3" "3.315.31"
Set-TextValue "E3" "  +0.66%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "592.15"
Set-TextValue "E5" "  +3.28%  "

Set-TextValue "D6" "181.83"
Set-TextValue "E6" "  +2.42%  "

Set-TextValue "D7" "0.638"
Set-TextValue "E7" "  +0.76%  "

Set-TextValue "E8" "  -0.03%  "

Set-TextValue "D9" "3.312.92"
Set-TextValue "E9" "  +0.71%  "

Set-TextValue "E10" "  +0.57%  "

Set-TextValue "E11" "  +3.12%  "

Set-TextValue "E12" "  +0.60%  "

Set-TextValue "D13" "3.888.76"
Set-TextValue "E13" "  +0.66%  "

Set-TextValue "E14" "  -2.20%  "

Set-TextValue "D15" "66.655.48"
Set-TextValue "E15" "  +1.12%  "

Set-TextValue "D16" "26.76"
Set-TextValue "E16" "  +0.40%  "

Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000164"
Set-TextValue "E17" "  +0.31%  "

Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.299.24"
Set-TextValue "E18" "  +1.00%  "

Set-TextValue "D19" "430.40"
Set-TextValue "E19" "  -1.45%  "

Set-TextValue "D20" "5.49"
Set-TextValue "E20" "  -1.76%  "

Set-TextValue "D21" "13.06"
Set-TextValue "E21" "  -1.76%  "

Set-TextValue "D22" "7.31"
Set-TextValue "E22" "  -1.59%  "

Set-TextValue "E23" "  +0.20%  "

Set-TextValue "D24" "71.51"
Set-TextValue "E24" "  -1.25%  "

Set-TextValue "E25" "  +1.03%  "

Set-TextValue "D26" "3.449.30"
Set-TextValue "E26" "  +0.13%  "

Set-TextValue "E27" "  +0.25%  "

Set-TextValue "D28" "0.207"
Set-TextValue "E28" "  +6.43%  "

Set-TextValue "E29" "  +0.88%  "

Set-TextValue "D30" "9.23"
Set-TextValue "E30" "  +3.24%  "

Set-TextValue "E31" "  -0.20%  "

Set-TextValue "D32" "1.93"
Set-TextValue "E32" "  -0.58%  "

Set-TextValue "D33" "22.40"
Set-TextValue "E33" "  +0.05%  "

Set-TextValue "E34" "  +0.10%  "

Set-TextValue "D35" "5.19"
Set-TextValue "E35" "  +1.07%  "

Set-TextValue "E36" "  -0.83%  "

Set-TextValue "D37" "1.19"
Set-TextValue "E37" "  -0.57%  "

Set-TextValue "D38" "158.93"
Set-TextValue "E38" "  +0.96%  "

Set-TextValue "D40" "2.872.35"
Set-TextValue "E40" "  +3.26%  "

Set-TextValue "E41" "  -0.15%  "

Set-TextValue "D42" "26.41"
Set-TextValue "E42" "  -2.12%  "

Set-TextValue "D43" "4.34"
Set-TextValue "E43" "  -0.39%  "

Set-TextValue "D44" "0.753"
Set-TextValue "E44" "  -3.86%  "

Set-TextValue "E45" "  -1.26%  "

Set-TextValue "D46" "5.93"
Set-TextValue "E46" "  -2.53%  "

Set-TextValue "B47" "dogwifhat"
Set-TextValue "C47" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D47" "2.32"
Set-TextValue "E47" "  +1.21%  "

Set-TextValue "B48" "Hedera"
Set-TextValue "C48" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D48" "0.0639"
Set-TextValue "E48" "  -3.17%  "

Set-TextValue "D49" "312.81"
Set-TextValue "E49" "  -2.47%  "

Set-TextValue "D50" "22.94"
Set-TextValue "E50" "  -2.48%  "

Set-TextValue "E51" "  -0.16%  "
